$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3076.6667
$ws.Range("J40").Value = 3153.3333
$ws.Range("L40").Value = 3153.3333
$ws.Range("N40").Value = -3503.3333

$ws.Range("H62").Value = 4882.9375
$ws.Range("I62").Value = 4882.9375
$ws.Range("K62").Value = 4882.9375
$ws.Range("M62").Value = -4258.9375

$ws.Range("H65").Value = 4882.9375
$ws.Range("I65").Value = 4882.9375
$ws.Range("K65").Value = 24414.6875
$ws.Range("M65").Value = -21294.6875

$ws.Range("H100").Value = 811.25
$ws.Range("I100").Value = 814.61536
$ws.Range("K100").Value = 814.61536
$ws.Range("M100").Value = -273.61536

$ws.Range("H112").Value = 1816.2667
$ws.Range("J112").Value = 1816.2667
$ws.Range("L112").Value = 5448.800099999999
$ws.Range("N112").Value = -7664.800099999999

$ws.Range("H135").Value = 711.3158
$ws.Range("I135").Value = 674.4
$ws.Range("J135").Value = 849.75
$ws.Range("K135").Value = 6069.599999999999
$ws.Range("L135").Value = 7647.75
$ws.Range("M135").Value = -3534.599999999999
$ws.Range("N135").Value = -12717.75

$ws.Range("H138").Value = 1958.7653
$ws.Range("I138").Value = 1341.4706
$ws.Range("J138").Value = 2286.7031
$ws.Range("K138").Value = 4024.4118
$ws.Range("L138").Value = 6860.1093
$ws.Range("M138").Value = 1115.5882
$ws.Range("N138").Value = -17140.1093

$ws.Range("H141").Value = 1752632
$ws.Range("I141").Value = 3501227.5
$ws.Range("J141").Value = 4036.5
$ws.Range("K141").Value = 10503682.5
$ws.Range("L141").Value = 12109.5
$ws.Range("M141").Value = -10498502.5
$ws.Range("N141").Value = -22469.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1031.4286
$ws.Range("J4").Value = 737.5
$ws.Range("L4").Value = 737.5
$ws.Range("N4").Value = -969.5

$ws.Range("H32").Value = 3422.09
$ws.Range("I32").Value = 3166.2632
$ws.Range("J32").Value = 8282.799999999999
$ws.Range("K32").Value = 3166.2632
$ws.Range("L32").Value = 8282.799999999999
$ws.Range("M32").Value = -2879.2632
$ws.Range("N32").Value = -8856.799999999999

$ws.Range("H61").Value = 34328.8
$ws.Range("I61").Value = 43253.684
$ws.Range("J61").Value = 6066.6665
$ws.Range("K61").Value = 43253.684
$ws.Range("L61").Value = 6066.6665
$ws.Range("M61").Value = -43041.684
$ws.Range("N61").Value = -6490.6665

$ws.Range("H97").Value = 1174.2609
$ws.Range("I97").Value = 949.3
$ws.Range("J97").Value = 2674
$ws.Range("K97").Value = 949.3
$ws.Range("L97").Value = 2674
$ws.Range("M97").Value = -453.3
$ws.Range("N97").Value = -3666

$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("M126").Value = -6530

$ws.Range("H136").Value = 34328.8
$ws.Range("I136").Value = 43253.684
$ws.Range("J136").Value = 6066.6665
$ws.Range("K136").Value = 129761.052
$ws.Range("L136").Value = 18199.9995
$ws.Range("M136").Value = -127211.052
$ws.Range("N136").Value = -23299.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3028.1904
$ws.Range("I20").Value = 2672.2
$ws.Range("J20").Value = 3918.1667
$ws.Range("K20").Value = 2672.2
$ws.Range("L20").Value = 3918.1667
$ws.Range("M20").Value = -2425.2
$ws.Range("N20").Value = -4412.1667

$ws.Range("H105").Value = 2489.3076
$ws.Range("I105").Value = 2364.1365
$ws.Range("K105").Value = 2364.1365
$ws.Range("M105").Value = -617.1365000000001

$ws.Range("H134").Value = 3426.527
$ws.Range("I134").Value = 3467.8948
$ws.Range("J134").Value = 3287.8235
$ws.Range("K134").Value = 10403.6844
$ws.Range("L134").Value = 9863.470499999999
$ws.Range("M134").Value = -7868.6844
$ws.Range("N134").Value = -14933.4705

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2249.5386
$ws.Range("J31").Value = 3953.25
$ws.Range("L31").Value = 3953.25
$ws.Range("N31").Value = -4543.25

$ws.Range("H34").Value = 2249.5386
$ws.Range("J34").Value = 3953.25
$ws.Range("L34").Value = 3953.25
$ws.Range("N34").Value = -4357.25

$ws.Range("H59").Value = 18600
$ws.Range("J59").Value = 18600
$ws.Range("L59").Value = 18600
$ws.Range("N59").Value = -20890

$ws.Range("H105").Value = 1646.1666
$ws.Range("I105").Value = 1646.1666
$ws.Range("K105").Value = 1646.1666
$ws.Range("M105").Value = 100.8334

$ws.Range("H134").Value = 1300.5714
$ws.Range("J134").Value = 1494.1818
$ws.Range("L134").Value = 4482.5454
$ws.Range("N134").Value = -9552.545399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 23818.75
$ws.Range("J9").Value = 23818.75
$ws.Range("L9").Value = 71456.25
$ws.Range("N9").Value = -71904.25

$ws.Range("H81").Value = 2743.8572
$ws.Range("J81").Value = 3001.1667
$ws.Range("L81").Value = 9003.500100000001
$ws.Range("N81").Value = -11249.5001

$ws.Range("H84").Value = 2743.8572
$ws.Range("J84").Value = 3001.1667
$ws.Range("L84").Value = 27010.5003
$ws.Range("N84").Value = -38242.5003

$ws.Range("H122").Value = 1022.5769
$ws.Range("J122").Value = 1073.4783
$ws.Range("L122").Value = 9661.304700000001
$ws.Range("N122").Value = -14561.3047

$ws.Range("H140").Value = 4380.5557
$ws.Range("I140").Value = 3696.0908
$ws.Range("K140").Value = 11088.2724
$ws.Range("M140").Value = -5908.2724

$ws.Range("H141").Value = 2719.2144
$ws.Range("I141").Value = 2466.8462
$ws.Range("K141").Value = 7400.5386
$ws.Range("M141").Value = -2220.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 4298707
$ws.Range("I11").Value = 5185709.5
$ws.Range("K11").Value = 5185709.5
$ws.Range("M11").Value = -5185570.5

$ws.Range("H80").Value = 1990
$ws.Range("I80").Value = 1990
$ws.Range("K80").Value = 1990
$ws.Range("M80").Value = -992

$ws.Range("H83").Value = 1990
$ws.Range("I83").Value = 1990
$ws.Range("K83").Value = 9950
$ws.Range("M83").Value = -4958

$ws.Range("H97").Value = 1689.3846
$ws.Range("I97").Value = 2153.1667
$ws.Range("J97").Value = 1291.8572
$ws.Range("K97").Value = 2153.1667
$ws.Range("L97").Value = 1291.8572
$ws.Range("M97").Value = -1657.1667
$ws.Range("N97").Value = -2283.8572

$ws.Range("H126").Value = 1952596.5
$ws.Range("I126").Value = 2224780
$ws.Range("K126").Value = 6674340
$ws.Range("M126").Value = -6671870

$ws.Range("H132").Value = 1427293.9
$ws.Range("I132").Value = 2405984
$ws.Range("K132").Value = 7217952
$ws.Range("M132").Value = -7215422

$ws.Range("H136").Value = 7360.826
$ws.Range("J136").Value = 7360.826
$ws.Range("L136").Value = 22082.478
$ws.Range("N136").Value = -27182.478

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 355797.8
$ws.Range("J2").Value = 91427.164
$ws.Range("L2").Value = 91427.164
$ws.Range("N2").Value = -91651.164

$ws.Range("H16").Value = 9623.429
$ws.Range("I16").Value = 13117.6
$ws.Range("K16").Value = 13117.6
$ws.Range("M16").Value = -12947.6

$ws.Range("H55").Value = 583.5789
$ws.Range("I55").Value = 493
$ws.Range("J55").Value = 779.8333
$ws.Range("K55").Value = 493
$ws.Range("L55").Value = 779.8333
$ws.Range("M55").Value = -320
$ws.Range("N55").Value = -1125.8333

$ws.Range("H82").Value = 1739.8572
$ws.Range("I82").Value = 1739.8572
$ws.Range("K82").Value = 1739.8572
$ws.Range("M82").Value = -1378.8572

$ws.Range("H85").Value = 1739.8572
$ws.Range("I85").Value = 1739.8572
$ws.Range("K85").Value = 1739.8572
$ws.Range("M85").Value = -491.8571999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 10872.111
$ws.Range("I96").Value = 3299.5
$ws.Range("J96").Value = 13035.714
$ws.Range("K96").Value = 3299.5
$ws.Range("L96").Value = 13035.714
$ws.Range("M96").Value = -1926.5
$ws.Range("N96").Value = -15781.714

$ws.Range("H122").Value = 55505.535
$ws.Range("I122").Value = 68756.586
$ws.Range("J122").Value = 2501.3333
$ws.Range("K122").Value = 206269.758
$ws.Range("L122").Value = 7503.999899999999
$ws.Range("M122").Value = -203819.758
$ws.Range("N122").Value = -12403.9999
